$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two title rows at the top of the sheet (rows 1 and 2),
# plus the blank row 3, so the header row (old row 4) becomes row 1.
$ws.Rows("1:3").Delete()

# Update the selected cell to match the target workbook state.
$ws.Range("I8").Select()
